# Repaired Excel Export and Repaied TP, FP, TN, FN calculation
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Only US-False Pos. Neg.")

# --- Update existing row 2 ---
$ws.Range("B2").Value = "gpt-4-turbo"
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 3
$ws.Range("G2").Value = "0.) Dummy Value Describtion`n1.) The redundancy here stems from the repetition of the testing action and the tested entity, expressed in slightly different ways, likely for emphasis or to accommodate different perspectives`n"
$ws.Range("I2").Value = "No text found in source"
$ws.Range("K2").Value = "No text found in source"

# --- Add the new row 3 values ---
# Column A holds a text value that looks numeric ("00"); format it as text
# first so Excel keeps it as the string "00" (matching row 2's A column)
# instead of converting it to the number 0.
$ws.Range("A3").NumberFormat = "@"
$ws.Range("A3").Value = "00"

$ws.Range("B3").Value = "gpt-4-turbo"
$ws.Range("C3").Value = "False Positive"
$ws.Range("D3").Value = "Benefit"
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 3
$ws.Range("G3").Value = "0.) The redundancy in these sentences lies in the repetition of the action ('login') and the destination ('webpage' and 'website'), which convey the same idea using slightly different wording.`n1.) Dummy Value Describtion`n"
$ws.Range("H3").Value = "0.) I can login into the webpage.`n1.) I can print a document`n"
$ws.Range("I3").Value = "No text found in source"
$ws.Range("J3").Value = "0.) I could login into the website`n1.) I can give the order to print`n"
$ws.Range("K3").Value = "No text found in source"

# --- Make row 3 share the exact same formatting as row 2 ---
$ws.Range("A2:K2").Copy()
$ws.Range("A3:K3").PasteSpecial(-4122)
$excel.CutCopyMode = 0

Write-Output "done"
